$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------
# 1. Data changes on "Planning & Journal" sheet
#    Row 10 ("Mise en place de Tailwind CSS"): now finished, half a
#    planned hour, half an hour actually spent.
# -------------------------------------------------------------------
$ws.Range("C10").Value = "Terminé"
$ws.Range("D10").Value = 0.5
$ws.Range("E10").Value = 0.5

#    Row 11 (Planification / Macro-planning) moves from "À faire" to
#    "En cours".
$ws.Range("C11").Value = "En cours"

# -------------------------------------------------------------------
# 2. Selection moves to H9 (cosmetic, matches the saved cursor state)
# -------------------------------------------------------------------
$ws.Range("H9").Select()

# -------------------------------------------------------------------
# 3. Conditional formatting rework on column C.
#    Previously C9 ("Terminé") shared its formatting rules with the
#    rest of the column; now C9 gets its own dedicated set of three
#    rules (Terminé / En cours / à faire), while the generic rules are
#    re-scoped to skip C9 (but now also cover C10, which used to have
#    its own separate "Terminé" rule).
# -------------------------------------------------------------------
$ws.Range("C1:C1048576").FormatConditions.Delete()

$greenFont = 24832
$greenFill = 13561798
$yellowFont = 22428
$yellowFill = 10284031
$redFont = 393372
$redFill = 13551615

# -- Rule: cellIs "Terminé" -> green, applies to C18:C1048576, C1:C8, C10:C15
foreach ($addr in @("C18:C1048576", "C1:C8", "C10:C15")) {
    $fc = $ws.Range($addr).FormatConditions.Add(1, 3, '="Terminé"')
    $fc.Priority = 11
    $fc.Font.Color = $greenFont
    $fc.Interior.Color = $greenFill
}

# -- Rule: containsText "En cours" -> yellow, applies to C17:C1048576, C1:C8, C10:C15
foreach ($addr in @("C17:C1048576", "C1:C8", "C10:C15")) {
    $fc = $ws.Range($addr).FormatConditions.Add(9, 0, "En cours")
    $fc.Priority = 9
    $fc.Text = "En cours"
    $anchor = ($addr -split ":")[0]
    $fc.Formula1 = '=NOT(ISERROR(SEARCH("En cours",' + $anchor + ')))'
    $fc.Font.Color = $yellowFont
    $fc.Interior.Color = $yellowFill
}

# -- Rule: containsText "à faire" -> red, applies to C1:C8, C10:C1048576
foreach ($addr in @("C1:C8", "C10:C1048576")) {
    $fc = $ws.Range($addr).FormatConditions.Add(9, 0, "à faire")
    $fc.Priority = 8
    $fc.Text = "à faire"
    $anchor = ($addr -split ":")[0]
    $fc.Formula1 = '=NOT(ISERROR(SEARCH("à faire",' + $anchor + ')))'
    $fc.Font.Color = $redFont
    $fc.Interior.Color = $redFill
}

# -- C9-only rules (Terminé / En cours / à faire)
$fc9a = $ws.Range("C9").FormatConditions.Add(1, 3, '="Terminé"')
$fc9a.Priority = 3
$fc9a.Font.Color = $greenFont
$fc9a.Interior.Color = $greenFill

$fc9b = $ws.Range("C9").FormatConditions.Add(9, 0, "En cours")
$fc9b.Priority = 2
$fc9b.Text = "En cours"
$fc9b.Formula1 = '=NOT(ISERROR(SEARCH("En cours",C9)))'
$fc9b.Font.Color = $yellowFont
$fc9b.Interior.Color = $yellowFill

$fc9c = $ws.Range("C9").FormatConditions.Add(9, 0, "à faire")
$fc9c.Priority = 1
$fc9c.Text = "à faire"
$fc9c.Formula1 = '=NOT(ISERROR(SEARCH("à faire",C9)))'
$fc9c.Font.Color = $redFont
$fc9c.Interior.Color = $redFill
